# Replace the blind vias with through hole ones
# Re-save touch-up: rename sheet, tidy column widths / row heights with wrap
# text on long descriptions, adjust zoom & selection, and nudge the logo
# picture to match the new column B width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab to match the board name.
$ws.Name = "H15R80"

# 2) Resize columns to the new layout.
$ws.Columns.Item(1).ColumnWidth = 11.21875
$ws.Columns.Item(2).ColumnWidth = 59.33203125
$ws.Columns.Item(3).ColumnWidth = 15.44140625
$ws.Columns.Item(4).ColumnWidth = 23.5546875
$ws.Columns.Item(5).ColumnWidth = 70.77734375
$ws.Columns.Item(6).ColumnWidth = 7.88671875
$ws.Columns.Item(7).ColumnWidth = 49

# 3) Wrap long text so the rows auto-size like in Excel after a resave.
$ws.Range("A9:G33").WrapText = $true
$ws.Range("A9:G33").VerticalAlignment = -4108  # xlCenter

# 4) Row heights to match the rewrapped content.
$ws.Rows.Item(2).RowHeight = 46.2
$ws.Rows.Item(8).RowHeight = 13.2
$ws.Rows.Item(10).RowHeight = 26.4
$ws.Rows.Item(11).RowHeight = 28.8
$ws.Rows.Item(12).RowHeight = 28.8
$ws.Rows.Item(13).RowHeight = 28.8
$ws.Rows.Item(16).RowHeight = 39.6
$ws.Rows.Item(17).RowHeight = 39.6
$ws.Rows.Item(18).RowHeight = 28.8
$ws.Rows.Item(19).RowHeight = 28.8
$ws.Rows.Item(20).RowHeight = 26.4
$ws.Rows.Item(24).RowHeight = 26.4
$ws.Rows.Item(25).RowHeight = 26.4
$ws.Rows.Item(26).RowHeight = 26.4
$ws.Rows.Item(27).RowHeight = 28.8
$ws.Rows.Item(28).RowHeight = 28.8
$ws.Rows.Item(29).RowHeight = 28.8
$ws.Rows.Item(30).RowHeight = 26.4
$ws.Rows.Item(31).RowHeight = 26.4

# 5) View tweaks.
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("B10").Select()

Write-Host "Applied H15R80 BOM touch-up edits"
